$wb = $excel.ActiveWorkbook

$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$boundariesSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "Version" banner and the "Recommended Citation" text on the About sheet.
foreach ($addr in @("A2", "A6")) {
    $cell = $aboutSheet.Range($addr)
    $current = $cell.Value2
    if ($current -ne $null -and $current.Contains($oldText)) {
        $cell.Value = $current.Replace($oldText, $newText)
    }
}

# Update the build_version column (S) for every data row on the Boundaries sheet.
$lastRow = $boundariesSheet.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $boundariesSheet.Cells.Item($row, 19)
    $current = $cell.Value2
    if ($current -ne $null -and $current.Contains($oldText)) {
        $cell.Value = $current.Replace($oldText, $newText)
    }
}
